$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-9 from 45204 to 45207
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
